# (#33) Alteração nos rótulos da tabela para já transformar a primeira linha
# em cabeçalho automaticamente no Power BI.
#
# For every worksheet, the header row (row 1) labels that are just a bare
# year ("2015", "2030", "2040", "2050") get prefixed with "Ano " (Year),
# while the worksheet whose header row uses date ranges ("2015-2030",
# "2031-2040", "2041-2050") gets its labels prefixed with "Intervalo "
# (Interval) instead. The first column's header (e.g. "Fonte/Tecnologia",
# "Período", "Tipo Expansão") is left untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastCol = $usedRange.Columns.Count

    # Decide which prefix this sheet's header row uses by inspecting the
    # text already in row 1 (column B onward): ranges use a hyphen.
    $prefix = "Ano "
    for ($c = 2; $c -le $lastCol; $c++) {
        $headerText = [string]$ws.Cells.Item(1, $c).Value2
        if ($headerText -match "-") {
            $prefix = "Intervalo "
            break
        }
    }

    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item(1, $c)
        $current = [string]$cell.Value2
        if ($current -notmatch "^Ano " -and $current -notmatch "^Intervalo ") {
            $cell.Value = "$prefix$current"
        }
    }
}
